$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.026.56"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.012.19"
$ws.Range("E3").Value = "  -1.85%  "

$ws.Range("E4").Value = "  -0.68%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.50"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.606"
$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.19"
$ws.Range("E8").Value = "  -1.83%  "

$ws.Range("E9").Value = "  -3.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0777"
$ws.Range("E10").Value = "  -3.44%  "

$ws.Range("E11").Value = "  -4.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.309.35"
$ws.Range("E12").Value = "  -1.84%  "

$ws.Range("E13").Value = "  -3.36%  "

$ws.Range("E14").Value = "  -3.91%  "

$ws.Range("E15").Value = "  -2.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.15"
$ws.Range("E16").Value = "  -1.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.014.38"
$ws.Range("E17").Value = "  -1.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.944.11"
$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  +2.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.85"
$ws.Range("E20").Value = "  -1.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0811"
$ws.Range("E21").Value = "  -3.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.75"
$ws.Range("E22").Value = "  -1.27%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  +2.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  -4.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.57"
$ws.Range("E26").Value = "  -2.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.96"
$ws.Range("E27").Value = "  -5.61%  "

$ws.Range("E28").Value = "  -2.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.56"
$ws.Range("E29").Value = "  -1.80%  "

$ws.Range("E30").Value = "  -4.79%  "

$ws.Range("E31").Value = "  -0.89%  "

$ws.Range("E32").Value = "  -2.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0600"
$ws.Range("E33").Value = "  -1.87%  "

$ws.Range("E34").Value = "  -1.84%  "

$ws.Range("E35").Value = "  -2.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.87"
$ws.Range("E36").Value = "  +2.24%  "

$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("E38").Value = "  -1.72%  "

$ws.Range("E39").Value = "  -0.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.463.26"
$ws.Range("E40").Value = "  -2.47%  "

$ws.Range("E41").Value = "  -4.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.34"
$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0911"
$ws.Range("E43").Value = "  -2.50%  "

$ws.Range("E44").Value = "  -4.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.15"
$ws.Range("E45").Value = "  +10.34%  "

$ws.Range("E48").Value = "  -1.71%  "

$ws.Range("E49").Value = "  -2.40%  "

$ws.Range("E50").Value = "  -1.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.197.60"

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.89"
$ws.Range("E46").Value = "  -5.43%  "

$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  -2.43%  "

